$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4.361249923706055
$ws.Range("E2").Value = 4.000625133514404
$ws.Range("F2").Value = 4.385313034057617
$ws.Range("G2").Value = 3.978125095367432
$ws.Range("H2").Value = 1256865381
$ws.Range("I2").Value = "ANET"

$ws.Range("D3").Value = 5.162499904632568
$ws.Range("E3").Value = 5.279375076293945
$ws.Range("F3").Value = 5.515625
$ws.Range("G3").Value = 4.788750171661377
$ws.Range("H3").Value = 1256865381
$ws.Range("I3").Value = "ANET"

$ws.Range("D4").Value = 3.819999933242798
$ws.Range("E4").Value = 4.031875133514404
$ws.Range("F4").Value = 4.3125
$ws.Range("G4").Value = 3.731250047683716
$ws.Range("H4").Value = 1256865381
$ws.Range("I4").Value = "ANET"

$ws.Range("D5").Value = 4.865624904632568
$ws.Range("E5").Value = 3.751874923706055
$ws.Range("F5").Value = 4.951562881469727
$ws.Range("G5").Value = 3.646874904632568
$ws.Range("H5").Value = 1256865381
$ws.Range("I5").Value = "ANET"

$ws.Range("D6").Value = 3.898124933242798
$ws.Range("E6").Value = 4.163750171661377
$ws.Range("F6").Value = 4.177499771118164
$ws.Range("G6").Value = 3.787499904632568
$ws.Range("H6").Value = 1256865381
$ws.Range("I6").Value = "ANET"

$ws.Range("D7").Value = 4.034999847412109
$ws.Range("E7").Value = 4.45437479019165
$ws.Range("F7").Value = 4.494375228881836
$ws.Range("G7").Value = 3.888124942779541
$ws.Range("H7").Value = 1256865381
$ws.Range("I7").Value = "ANET"

$ws.Range("D8").Value = 5.306875228881836
$ws.Range("E8").Value = 5.296875
$ws.Range("F8").Value = 5.476250171661377
$ws.Range("G8").Value = 4.926249980926514
$ws.Range("H8").Value = 1256865381
$ws.Range("I8").Value = "ANET"

$ws.Range("D9").Value = 6.096250057220459
$ws.Range("E9").Value = 5.875
$ws.Range("F9").Value = 6.4375
$ws.Range("G9").Value = 5.458125114440918
$ws.Range("H9").Value = 1256865381
$ws.Range("I9").Value = "ANET"

$ws.Range("D10").Value = 8.233124732971191
$ws.Range("E10").Value = 8.727499961853027
$ws.Range("F10").Value = 8.770000457763672
$ws.Range("G10").Value = 8.041250228881836
$ws.Range("H10").Value = 1256865381
$ws.Range("I10").Value = "ANET"

$ws.Range("D11").Value = 9.441249847412109
$ws.Range("E11").Value = 9.330624580383301
$ws.Range("F11").Value = 9.9375
$ws.Range("G11").Value = 9.011875152587891
$ws.Range("H11").Value = 1256865381
$ws.Range("I11").Value = "ANET"

$ws.Range("D12").Value = 11.875
$ws.Range("E12").Value = 12.49312496185303
$ws.Range("F12").Value = 12.62187480926514
$ws.Range("G12").Value = 11.578125
$ws.Range("H12").Value = 1256865381
$ws.Range("I12").Value = "ANET"

$ws.Range("D13").Value = 14.72375011444092
$ws.Range("E13").Value = 17.23875045776367
$ws.Range("F13").Value = 17.79750061035156
$ws.Range("G13").Value = 14.17249965667725
$ws.Range("H13").Value = 1256865381
$ws.Range("I13").Value = "ANET"

$ws.Range("D14").Value = 15.8412504196167
$ws.Range("E14").Value = 16.53437423706055
$ws.Range("F14").Value = 17.00625038146973
$ws.Range("G14").Value = 14.91874980926514
$ws.Range("H14").Value = 1256865381
$ws.Range("I14").Value = "ANET"

$ws.Range("D15").Value = 15.93624973297119
$ws.Range("E15").Value = 15.98312473297119
$ws.Range("F15").Value = 17.57500076293945
$ws.Range("G15").Value = 15.765625
$ws.Range("H15").Value = 1256865381
$ws.Range("I15").Value = "ANET"

$ws.Range("D16").Value = 16.70062446594238
$ws.Range("E16").Value = 14.39687538146973
$ws.Range("F16").Value = 16.80125045776367
$ws.Range("G16").Value = 13.31937503814697
$ws.Range("H16").Value = 1256865381
$ws.Range("I16").Value = "ANET"

$ws.Range("D17").Value = 12.8125
$ws.Range("E17").Value = 13.42374992370606
$ws.Range("F17").Value = 14.57250022888184
$ws.Range("G17").Value = 12.45031261444092
$ws.Range("H17").Value = 1256865381
$ws.Range("I17").Value = "ANET"

$ws.Range("D18").Value = 19.92124938964844
$ws.Range("E18").Value = 19.51812553405762
$ws.Range("F18").Value = 20.70437431335449
$ws.Range("G18").Value = 19.15393829345703
$ws.Range("H18").Value = 1256865381
$ws.Range("I18").Value = "ANET"

$ws.Range("D19").Value = 16.59187507629395
$ws.Range("E19").Value = 17.09062576293945
$ws.Range("F19").Value = 18.09062576293945
$ws.Range("G19").Value = 16.33437538146973
$ws.Range("H19").Value = 1256865381
$ws.Range("I19").Value = "ANET"

$ws.Range("D20").Value = 15.01687526702881
$ws.Range("E20").Value = 15.28562545776367
$ws.Range("F20").Value = 15.73375034332275
$ws.Range("G20").Value = 13.45374965667725
$ws.Range("H20").Value = 1256865381
$ws.Range("I20").Value = "ANET"

$ws.Range("D21").Value = 12.85499954223633
$ws.Range("E21").Value = 13.95874977111816
$ws.Range("F21").Value = 15.11812496185303
$ws.Range("G21").Value = 12.375
$ws.Range("H21").Value = 1256865381
$ws.Range("I21").Value = "ANET"

$ws.Range("D22").Value = 12.16187477111816
$ws.Range("E22").Value = 13.70625019073486
$ws.Range("F22").Value = 14
$ws.Range("G22").Value = 11.81187534332275
$ws.Range("H22").Value = 1256865381
$ws.Range("I22").Value = "ANET"

$ws.Range("D23").Value = 13.16187477111816
$ws.Range("E23").Value = 16.23562431335449
$ws.Range("F23").Value = 16.23749923706055
$ws.Range("G23").Value = 12.9556245803833
$ws.Range("H23").Value = 1256865381
$ws.Range("I23").Value = "ANET"

$ws.Range("D24").Value = 13.03125
$ws.Range("E24").Value = 13.05624961853027
$ws.Range("F24").Value = 14.57999992370606
$ws.Range("G24").Value = 12.52187538146973
$ws.Range("H24").Value = 1256865381
$ws.Range("I24").Value = "ANET"

$ws.Range("D25").Value = 18.17250061035156
$ws.Range("E25").Value = 19.22249984741211
$ws.Range("F25").Value = 20.02124977111816
$ws.Range("G25").Value = 17.52687454223633
$ws.Range("H25").Value = 1256865381
$ws.Range("I25").Value = "ANET"

$ws.Range("D26").Value = 19.00374984741211
$ws.Range("E26").Value = 19.6981258392334
$ws.Range("F26").Value = 20.21187591552734
$ws.Range("G26").Value = 18.9637508392334
$ws.Range("H26").Value = 1256865381
$ws.Range("I26").Value = "ANET"

$ws.Range("D27").Value = 22.71187591552734
$ws.Range("E27").Value = 23.77437591552734
$ws.Range("F27").Value = 23.9715633392334
$ws.Range("G27").Value = 22.33187484741211
$ws.Range("H27").Value = 1256865381
$ws.Range("I27").Value = "ANET"

$ws.Range("D28").Value = 21.5625
$ws.Range("E28").Value = 25.60562515258789
$ws.Range("F28").Value = 25.69874954223633
$ws.Range("G28").Value = 21.29500007629395
$ws.Range("H28").Value = 1256865381
$ws.Range("I28").Value = "ANET"

$ws.Range("D29").Value = 36.02249908447266
$ws.Range("E29").Value = 31.07749938964844
$ws.Range("F29").Value = 36.125
$ws.Range("G29").Value = 28.71249961853028
$ws.Range("H29").Value = 1256865381
$ws.Range("I29").Value = "ANET"

$ws.Range("D30").Value = 35.08250045776367
$ws.Range("E30").Value = 28.89249992370605
$ws.Range("F30").Value = 35.89250183105469
$ws.Range("G30").Value = 28.34749984741211
$ws.Range("H30").Value = 1256865381
$ws.Range("I30").Value = "ANET"

$ws.Range("D31").Value = 23.28499984741211
$ws.Range("E31").Value = 29.15749931335449
$ws.Range("F31").Value = 29.33250045776367
$ws.Range("G31").Value = 22.82749938964844
$ws.Range("H31").Value = 1256865381
$ws.Range("I31").Value = "ANET"

$ws.Range("D32").Value = 28.71750068664551
$ws.Range("E32").Value = 30.21500015258789
$ws.Range("F32").Value = 31.63249969482422
$ws.Range("G32").Value = 24.55125045776367
$ws.Range("H32").Value = 1256865381
$ws.Range("I32").Value = "ANET"

$ws.Range("D33").Value = 30.72750091552734
$ws.Range("E33").Value = 31.5049991607666
$ws.Range("F33").Value = 31.76749992370605
$ws.Range("G33").Value = 27.22750091552734
$ws.Range("H33").Value = 1256865381
$ws.Range("I33").Value = "ANET"

$ws.Range("D34").Value = 41.41249847412109
$ws.Range("E34").Value = 40.04000091552734
$ws.Range("F34").Value = 42.16500091552734
$ws.Range("G34").Value = 38.08750152587891
$ws.Range("H34").Value = 1256865381
$ws.Range("I34").Value = "ANET"

$ws.Range("D35").Value = 40.33499908447266
$ws.Range("E35").Value = 38.77249908447266
$ws.Range("F35").Value = 44.21125030517578
$ws.Range("G35").Value = 37.27500152587891
$ws.Range("H35").Value = 1256865381
$ws.Range("I35").Value = "ANET"

$ws.Range("D36").Value = 46.13000106811523
$ws.Range("E36").Value = 50.09249877929688
$ws.Range("F36").Value = 50.36249923706055
$ws.Range("G36").Value = 42.0625
$ws.Range("H36").Value = 1256865381
$ws.Range("I36").Value = "ANET"

$ws.Range("D37").Value = 58.5
$ws.Range("E37").Value = 64.66999816894531
$ws.Range("F37").Value = 68.22499847412109
$ws.Range("G37").Value = 56.65000152587891
$ws.Range("H37").Value = 1256865381
$ws.Range("I37").Value = "ANET"

$ws.Range("D38").Value = 73.06999969482422
$ws.Range("E38").Value = 64.13999938964844
$ws.Range("F38").Value = 75.875
$ws.Range("G38").Value = 60.08000183105469
$ws.Range("H38").Value = 1256865381
$ws.Range("I38").Value = "ANET"

$ws.Range("D39").Value = 88.67749786376953
$ws.Range("E39").Value = 86.63749694824219
$ws.Range("F39").Value = 94.125
$ws.Range("G39").Value = 76.77999877929688
$ws.Range("H39").Value = 1256865381
$ws.Range("I39").Value = "ANET"

$ws.Range("D40").Value = 97.09249877929688
$ws.Range("E40").Value = 96.61000061035156
$ws.Range("F40").Value = 105.682502746582
$ws.Range("G40").Value = 94.61499786376952
$ws.Range("H40").Value = 1256865381
$ws.Range("I40").Value = "ANET"

$ws.Range("D41").Value = 111.3199996948242
$ws.Range("E41").Value = 115.2300033569336
$ws.Range("F41").Value = 133.5749969482422
$ws.Range("G41").Value = 97.68000030517578
$ws.Range("H41").Value = 1256865381
$ws.Range("I41").Value = "ANET"

$ws.Range("D42").Value = 77.15000152587891
$ws.Range("E42").Value = 82.26999664306641
$ws.Range("F42").Value = 82.55500030517578
$ws.Range("G42").Value = 59.43000030517578
$ws.Range("H42").Value = 1256865381
$ws.Range("I42").Value = "ANET"

$ws.Range("D43").Value = 102.3099975585938
$ws.Range("E43").Value = 123.2200012207031
$ws.Range("F43").Value = 125.8099975585938
$ws.Range("G43").Value = 97.13999938964844
$ws.Range("H43").Value = 1256865381
$ws.Range("I43").Value = "ANET"

Write-Output "done"
